# Adds a new row of data (row 32) for a new timeline entry, using the same
# date as the previous row (31/3/2024), 5 hours spent, and a new description
# string. Rows 32-35 remain otherwise empty (as before), and rows 36-38
# (totals) are left where they are - only the SUM formula result changes
# because it now includes C32 in its range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting from the row above so the new row matches the style
# of surrounding data rows (without inserting/shifting any rows).
$ws.Range("B31:D31").Copy()
$ws.Range("B32:D32").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B32").Value = $ws.Range("B31").Value2
$ws.Range("C32").Value = 5
$ws.Range("D32").Value = "Added Exception handling with custom error codes for better FE control over errors"

# Update the current view/selection to match the new state
# (scrolled down so row 4 is at the top, with G29 as the active cell).
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G29").Select()
